# Insert a new "enable_year" configuration row right after the existing
# "enable_capacity" row (row 6), pushing the whole table down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 7 - everything currently at/after row 7
# shifts down by one.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 following the same pattern as the
# other "configuration" rows (e.g. row 6: enable_capacity).
$ws.Cells.Item(7, 1).Value = "CHE"
$ws.Cells.Item(7, 2).Value = "conv_transmission_elec"
$ws.Cells.Item(7, 3).Value = "enable_year"
$ws.Cells.Item(7, 4).Value = "configuration"
$ws.Cells.Item(7, 7).Value = 1990

# Move the selection, matching the post-edit cursor location.
$ws.Range("F12").Select()

$wb.Save()
